# Update 'F' column (想去人数 / interested count) values across sheets
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 295  # was 294
$ws.Range("F3").Value = 91  # was 90
$ws.Range("F4").Value = 1211  # was 1208
$ws.Range("F5").Value = 846  # was 841
$ws.Range("F6").Value = 871  # was 870
$ws.Range("F7").Value = 1591  # was 1585
$ws.Range("F8").Value = 331  # was 325
$ws.Range("F9").Value = 1087  # was 1081
$ws.Range("F12").Value = 219  # was 215
$ws.Range("F14").Value = 551  # was 546
$ws.Range("F15").Value = 92  # was 89
$ws.Range("F16").Value = 56  # was 55
$ws.Range("F20").Value = 600  # was 597
$ws.Range("F21").Value = 593  # was 592
$ws.Range("F22").Value = 79  # was 78
$ws.Range("F26").Value = 210  # was 207

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1054  # was 1052
$ws.Range("F4").Value = 295  # was 291
$ws.Range("F5").Value = 19  # was 18
$ws.Range("F6").Value = 211  # was 210
$ws.Range("F8").Value = 606  # was 604

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 277  # was 276

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 277  # was 276
$ws.Range("F3").Value = 295  # was 294
$ws.Range("F4").Value = 1054  # was 1052
$ws.Range("F5").Value = 91  # was 90
$ws.Range("F6").Value = 1211  # was 1208
$ws.Range("F7").Value = 846  # was 841
$ws.Range("F8").Value = 871  # was 870
$ws.Range("F9").Value = 1591  # was 1585
$ws.Range("F10").Value = 331  # was 325
$ws.Range("F11").Value = 1087  # was 1081
$ws.Range("F14").Value = 219  # was 215
$ws.Range("F16").Value = 551  # was 546
$ws.Range("F17").Value = 92  # was 89
$ws.Range("F18").Value = 56  # was 55
$ws.Range("F21").Value = 295  # was 291
$ws.Range("F24").Value = 19  # was 18
$ws.Range("F25").Value = 211  # was 210
$ws.Range("F26").Value = 211  # was 210
$ws.Range("F27").Value = 600  # was 597
$ws.Range("F28").Value = 593  # was 592
$ws.Range("F29").Value = 79  # was 78
$ws.Range("F34").Value = 210  # was 207
$ws.Range("F35").Value = 606  # was 604
